$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 600.3333
$ws.Range("I8").Value = 600.3333
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1800.9999
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1661.9999
$ws.Range("N8").Value = $null
$ws.Range("H58").Value = 60748.41
$ws.Range("J58").Value = 73757.36
$ws.Range("L58").Value = 221272.08
$ws.Range("N58").Value = -221572.08
$ws.Range("H74").Value = 4927.273
$ws.Range("I74").Value = 4414.2856
$ws.Range("J74").Value = 5825
$ws.Range("K74").Value = 4414.2856
$ws.Range("L74").Value = 5825
$ws.Range("M74").Value = -3478.2856
$ws.Range("N74").Value = -7697
$ws.Range("H77").Value = 4927.273
$ws.Range("I77").Value = 4414.2856
$ws.Range("J77").Value = 5825
$ws.Range("K77").Value = 22071.428
$ws.Range("L77").Value = 29125
$ws.Range("M77").Value = -17391.428
$ws.Range("N77").Value = -38485
$ws.Range("H135").Value = 57694550
$ws.Range("I135").Value = 25001970
$ws.Range("K135").Value = 225017730
$ws.Range("M135").Value = -225015195
$ws.Range("H137").Value = 4469.6333
$ws.Range("I137").Value = 2433.7334
$ws.Range("K137").Value = 7301.2002
$ws.Range("M137").Value = -4751.2002
$ws.Range("H138").Value = 3388.3333
$ws.Range("I138").Value = 1914.6923
$ws.Range("J138").Value = 3935.6858
$ws.Range("K138").Value = 5744.0769
$ws.Range("L138").Value = 11807.0574
$ws.Range("M138").Value = -604.0769
$ws.Range("N138").Value = -22087.0574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10814.667
$ws.Range("J3").Value = 10897.1
$ws.Range("L3").Value = 10897.1
$ws.Range("N3").Value = -11127.1
$ws.Range("H32").Value = 6858.108
$ws.Range("I32").Value = 6139.75
$ws.Range("K32").Value = 6139.75
$ws.Range("M32").Value = -5852.75
$ws.Range("H61").Value = 8314.638999999999
$ws.Range("I61").Value = 4376.3706
$ws.Range("J61").Value = 13631.3
$ws.Range("K61").Value = 4376.3706
$ws.Range("L61").Value = 13631.3
$ws.Range("M61").Value = -4164.3706
$ws.Range("N61").Value = -14055.3
$ws.Range("H63").Value = 6550.25
$ws.Range("I63").Value = 2733.6667
$ws.Range("K63").Value = 2733.6667
$ws.Range("M63").Value = -2047.6667
$ws.Range("H66").Value = 6550.25
$ws.Range("I66").Value = 2733.6667
$ws.Range("K66").Value = 13668.3335
$ws.Range("M66").Value = -10236.3335
$ws.Range("H74").Value = 57819.832
$ws.Range("I74").Value = 72457.5
$ws.Range("J74").Value = 9027.6
$ws.Range("K74").Value = 72457.5
$ws.Range("L74").Value = 9027.6
$ws.Range("M74").Value = -71583.5
$ws.Range("N74").Value = -10775.6
$ws.Range("H77").Value = 57819.832
$ws.Range("I77").Value = 72457.5
$ws.Range("J77").Value = 9027.6
$ws.Range("K77").Value = 362287.5
$ws.Range("L77").Value = 45138
$ws.Range("M77").Value = -357919.5
$ws.Range("N77").Value = -53874
$ws.Range("H132").Value = 6192.5835
$ws.Range("I132").Value = 2399.0715
$ws.Range("J132").Value = 8606.637000000001
$ws.Range("K132").Value = 7197.2145
$ws.Range("L132").Value = 25819.911
$ws.Range("M132").Value = -4667.2145
$ws.Range("N132").Value = -30879.911
$ws.Range("H136").Value = 8314.638999999999
$ws.Range("I136").Value = 4376.3706
$ws.Range("J136").Value = 13631.3
$ws.Range("K136").Value = 13129.1118
$ws.Range("L136").Value = 40893.89999999999
$ws.Range("M136").Value = -10579.1118
$ws.Range("N136").Value = -45993.89999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1008.75
$ws.Range("J12").Value = 1450
$ws.Range("L12").Value = 1450
$ws.Range("N12").Value = -1786
$ws.Range("H134").Value = 41373.117
$ws.Range("I134").Value = 3042.95
$ws.Range("J134").Value = 169140.33
$ws.Range("K134").Value = 9128.849999999999
$ws.Range("L134").Value = 507420.99
$ws.Range("M134").Value = -6593.849999999999
$ws.Range("N134").Value = -512490.99

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 100006
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 100006
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 100006
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = -100284
$ws.Range("H58").Value = 1785074.1
$ws.Range("I58").Value = 3369147.5
$ws.Range("J58").Value = 2991.6667
$ws.Range("K58").Value = 3369147.5
$ws.Range("L58").Value = 2991.6667
$ws.Range("M58").Value = -3368944.5
$ws.Range("N58").Value = -3397.6667
$ws.Range("H107").Value = 1144.909
$ws.Range("J107").Value = 890.75
$ws.Range("L107").Value = 890.75
$ws.Range("N107").Value = -4730.75
$ws.Range("H132").Value = 1719.7869
$ws.Range("I132").Value = 1564.5807
$ws.Range("J132").Value = 1880.1666
$ws.Range("K132").Value = 4693.742099999999
$ws.Range("L132").Value = 5640.4998
$ws.Range("M132").Value = -2163.742099999999
$ws.Range("N132").Value = -10700.4998
$ws.Range("H134").Value = 22137.037
$ws.Range("I134").Value = 41874.52
$ws.Range("J134").Value = 4514.2856
$ws.Range("K134").Value = 125623.56
$ws.Range("L134").Value = 13542.8568
$ws.Range("M134").Value = -123088.56
$ws.Range("N134").Value = -18612.8568
$ws.Range("H136").Value = 1785074.1
$ws.Range("I136").Value = 3369147.5
$ws.Range("J136").Value = 2991.6667
$ws.Range("K136").Value = 10107442.5
$ws.Range("L136").Value = 8975.000100000001
$ws.Range("M136").Value = -10104892.5
$ws.Range("N136").Value = -14075.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2293.2942
$ws.Range("I109").Value = 399
$ws.Range("J109").Value = 2876.1538
$ws.Range("K109").Value = 1197
$ws.Range("L109").Value = 8628.4614
$ws.Range("M109").Value = -157
$ws.Range("N109").Value = -10708.4614
$ws.Range("H131").Value = 14455.029
$ws.Range("I131").Value = 350.17242
$ws.Range("J131").Value = 96263.2
$ws.Range("K131").Value = 1050.51726
$ws.Range("L131").Value = 288789.6
$ws.Range("M131").Value = 3989.48274
$ws.Range("N131").Value = -298869.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5516.18
$ws.Range("I70").Value = 4666.6665
$ws.Range("J70").Value = 5570.4043
$ws.Range("K70").Value = 4666.6665
$ws.Range("L70").Value = 5570.4043
$ws.Range("M70").Value = -4396.6665
$ws.Range("N70").Value = -6110.4043
$ws.Range("H73").Value = 5516.18
$ws.Range("I73").Value = 4666.6665
$ws.Range("J73").Value = 5570.4043
$ws.Range("K73").Value = 4666.6665
$ws.Range("L73").Value = 5570.4043
$ws.Range("M73").Value = -3730.6665
$ws.Range("N73").Value = -7442.4043
$ws.Range("H122").Value = 10528.571
$ws.Range("I122").Value = 50000
$ws.Range("J122").Value = 3950
$ws.Range("K122").Value = 150000
$ws.Range("L122").Value = 11850
$ws.Range("M122").Value = -147550
$ws.Range("N122").Value = -16750
$ws.Range("H126").Value = 3090.85
$ws.Range("I126").Value = 2025
$ws.Range("J126").Value = 3801.4167
$ws.Range("K126").Value = 6075
$ws.Range("L126").Value = 11404.2501
$ws.Range("M126").Value = -3605
$ws.Range("N126").Value = -16344.2501
$ws.Range("H132").Value = 141073.88
$ws.Range("I132").Value = 203715.8
$ws.Range("J132").Value = 36670.668
$ws.Range("K132").Value = 611147.3999999999
$ws.Range("L132").Value = 110012.004
$ws.Range("M132").Value = -608617.3999999999
$ws.Range("N132").Value = -115072.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180
$ws.Range("H136").Value = 5939.457
$ws.Range("I136").Value = 4288.7896
$ws.Range("J136").Value = 7899.625
$ws.Range("K136").Value = 12866.3688
$ws.Range("L136").Value = 23698.875
$ws.Range("M136").Value = -10316.3688
$ws.Range("N136").Value = -28798.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4497.2085
$ws.Range("I132").Value = 4047.0715
$ws.Range("J132").Value = 5127.4
$ws.Range("K132").Value = 12141.2145
$ws.Range("L132").Value = 15382.2
$ws.Range("M132").Value = -9611.2145
$ws.Range("N132").Value = -20442.2
$ws.Range("H136").Value = 3959.116
$ws.Range("I136").Value = 1444.7556
$ws.Range("K136").Value = 4334.266799999999
$ws.Range("M136").Value = -1784.266799999999
